# "lots of little edits"
#  - bump the cached datetimeFigureOut text from 5/7/17 -> 8/30/17
#    everywhere it is cached (slide master(s), slide layouts, notes master)
#  - simplify the small "prog"/"stmt_list" diagram on slide 1:
#      * TextBox 5  (id=6)  "stmt_list" -> "prog", shrink its box
#      * TextBox 39 (id=40) "s"+"tmt"   -> "stmt" (merge runs)
#      * TextBox 40 (id=41) "stmt_list" -> "prog", shrink its box
#      * reroute/resize a few connectors + the quote textbox
#      * delete the standalone "prog" textbox (id=29) and the
#        connector that fed it (id=3)
#
# NOTE: this COM host doesn't give PS functions their own variable
# scope, so every loop below uses its own uniquely-named counter
# instead of nesting calls that would otherwise clobber an outer
# loop's counter.

$p = $ppt.ActivePresentation

# ---- 1. re-cache the date placeholders everywhere they live -------------

for ($designIx = 1; $designIx -le $p.Designs.Count; $designIx++) {
    $curDesign = $p.Designs.Item($designIx)
    $curMaster = $curDesign.SlideMaster

    for ($masterShapeIx = 1; $masterShapeIx -le $curMaster.Shapes.Count; $masterShapeIx++) {
        $masterShape = $curMaster.Shapes.Item($masterShapeIx)
        if ($masterShape.HasTextFrame) {
            if ($masterShape.TextFrame.TextRange.Text -eq "5/7/17") {
                $masterShape.TextFrame.TextRange.Text = "8/30/17"
            }
        }
    }

    for ($layoutIx = 1; $layoutIx -le $curMaster.CustomLayouts.Count; $layoutIx++) {
        $curLayout = $curMaster.CustomLayouts.Item($layoutIx)
        for ($layoutShapeIx = 1; $layoutShapeIx -le $curLayout.Shapes.Count; $layoutShapeIx++) {
            $layoutShape = $curLayout.Shapes.Item($layoutShapeIx)
            if ($layoutShape.HasTextFrame) {
                if ($layoutShape.TextFrame.TextRange.Text -eq "5/7/17") {
                    $layoutShape.TextFrame.TextRange.Text = "8/30/17"
                }
            }
        }
    }
}

$notesMaster = $p.NotesMaster
for ($notesShapeIx = 1; $notesShapeIx -le $notesMaster.Shapes.Count; $notesShapeIx++) {
    $notesShape = $notesMaster.Shapes.Item($notesShapeIx)
    if ($notesShape.HasTextFrame) {
        if ($notesShape.TextFrame.TextRange.Text -eq "5/7/17") {
            $notesShape.TextFrame.TextRange.Text = "8/30/17"
        }
    }
}

# ---- 2. slide 1 diagram tidy-up ------------------------------------------

$s = $p.Slides.Item(1)
$emu = 12700.0

$shape6  = $s.Shapes.Item("TextBox 5")
$shape40 = $s.Shapes.Item("TextBox 39")
$shape41 = $s.Shapes.Item("TextBox 40")
$conn45  = $s.Shapes.Item("Straight Connector 45")
$conn46  = $s.Shapes.Item("Straight Connector 46")
$tb104   = $s.Shapes.Item("TextBox 104")
$conn106 = $s.Shapes.Item("Straight Connector 106")
$tb28    = $s.Shapes.Item("TextBox 28")
$conn2   = $s.Shapes.Item("Straight Connector 2")

# TextBox 5 (id=6): "stmt_list" -> "prog", narrower box
$shape6.TextFrame.TextRange.Text = "prog"
$shape6.Width = 686068 / $emu

# TextBox 39 (id=40): merge "s" + "tmt" -> "stmt" (keep trailing space run)
$shape40.TextFrame.TextRange.Characters(1, 4).Text = "stmt"

# TextBox 40 (id=41): "stmt_list" -> "prog", narrower box
$shape41.TextFrame.TextRange.Text = "prog"
$shape41.Width = 686068 / $emu

# Straight Connector 45 (id=46): connects shape6 -> shape40, shrinks with shape6
$conn45.Width = 1308653 / $emu

# Straight Connector 46 (id=47): connects shape6 -> shape41, shifts with shape41
$conn46.Left = 4699137 / $emu

# TextBox 104 (id=105): the quote box moves left
$tb104.Left = 5817214 / $emu

# Straight Connector 106 (id=107): connects shape41 -> tb104, reroutes
$conn106.HorizontalFlip = $false
$conn106.Left = 6022586 / $emu
$conn106.Width = 25483 / $emu

# Drop the now-redundant standalone "prog" textbox (id=29) and the
# connector that used to link it down into shape6 (id=3).
$conn2.Delete()
$tb28.Delete()
